# Add a new row 5 to the "Client Data" sheet: a duplicate of row 4's client
# record, but with an updated "Date of Birth" value ("2025-03-30" instead of
# "2025-04-09").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 4 into row 5. A plain Range-to-Range Copy preserves the
# shared-string references, numeric cells and cell styling (style 0) exactly
# as they are on row 4, instead of re-typing every value (which would risk
# Excel "smart" reinterpretation of things like the phone number / numeric
# text).
$ws.Range("A4:J4").Copy($ws.Range("A5:J5"))

# Now fix up column B (Date of Birth) on the new row so it reads
# "2025-03-30". Assigning that string straight to a cell's .Value makes
# Excel auto-detect it as a date and store a date serial number with a date
# number format instead of the literal text - not what we want here. Using
# a formula that evaluates to the text avoids the date auto-detection, and
# copying that formula's *value only* (PasteSpecial xlPasteValues = -4163)
# into the target cell converts it to a plain literal string cell while
# keeping the destination's existing (default) style untouched.
$ws.Range("L1").Formula = "=""2025-03-30"""
$ws.Range("L1").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("L1").Clear()
